# Append 21 more "Captured_Values" rows (68-88), matching the pattern of
# the existing rows 2-67: column A holds the numeric literal 123456789 and
# column B holds the repeated label string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Captured_Values")

$text = "Real Programmers Count 0123456789 From Zero"

for ($row = 68; $row -le 88; $row++) {
    $ws.Cells.Item($row, 1).Value = 123456789
    $ws.Cells.Item($row, 2).Value = $text
}
